# The workbook contains a weekly price log for Mango at "Vega Central
# Mapocho de Santiago". A new observation (week) was inserted into the
# table right before the existing row 231, which pushes every row from
# 231 onward down by one (231->232, 232->233, ..., 253->254). The new
# row is populated with its own data while all the columns that are
# constant across the whole sheet (Mercado ID, Mercado, Region, Codreg,
# Tipo, Producto ID, Producto, Categoria ID, Categoria, Variedad, Unidad
# de comercializacion, Kg/unidad) are copied from the neighbouring row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 231; this shifts old rows 231-253
# down to 232-254 and extends the sheet dimension to A1:T254.
$ws.Rows.Item(231).Insert()

# Fill in the constant columns by copying them from row 232 (which now
# holds what used to be row 231's data, i.e. the same template values).
$ws.Range("A231").Value = $ws.Range("A232").Value()
$ws.Range("B231").Value = $ws.Range("B232").Value()
$ws.Range("C231").Value = $ws.Range("C232").Value()
$ws.Range("E231").Value = $ws.Range("E232").Value()
$ws.Range("F231").Value = $ws.Range("F232").Value()
$ws.Range("G231").Value = $ws.Range("G232").Value()
$ws.Range("H231").Value = $ws.Range("H232").Value()
$ws.Range("I231").Value = $ws.Range("I232").Value()
$ws.Range("J231").Value = $ws.Range("J232").Value()
$ws.Range("K231").Value = $ws.Range("K232").Value()
$ws.Range("Q231").Value = $ws.Range("Q232").Value()
$ws.Range("T231").Value = $ws.Range("T232").Value()

# New row's own data (date, quality, volume, min/max/avg price, origin,
# price per kg).
$ws.Range("D231").Value = 44461
$ws.Range("L231").Value = "Primera"
$ws.Range("M231").Value = 430
$ws.Range("N231").Value = 7500
$ws.Range("O231").Value = 8000
$ws.Range("P231").Value = 7709
$ws.Range("R231").Value = "Brasil"
$ws.Range("S231").Value = 1927
